$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of timesheet data (row 6)
$ws.Range("B6").Value = 1.5
$ws.Range("C6").Value = "11pm"
$ws.Range("D6").Value = "12:30pm"

# Move selection to reflect where the user ended up after entering the row
$ws.Range("D7").Select()
